# "maqueta de la pantalla de consumo"
# Populate the "consumo" (fuel-consumption screen) task list on Hoja2
# with the new backlog items, tweak a couple of existing progress
# values, and restore the expected selection/view state.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Hoja2: row 13 progress resets to 0 (real / est columns) ---------------
$ws2.Range("B13").Value = 0
$ws2.Range("C13").Value = 0

# --- Hoja2: new task rows 16-26 ---------------------------------------------
# B = real days, C = estimated days, D = task description (copy the D-column
# "task" formatting - style index 7 - from an existing styled cell so we
# don't mint a new duplicate style entry).

$tasks = @(
    @{ Row = 16; B = 0.2; C = 0.2; Text = "al seleccionar vehiculo, rellenar la caja rendimiento" },
    @{ Row = 17; B = 0.4; C = 0.4; Text = "al modificar distancia o rendimieno actualizar consumo en litros" },
    @{ Row = 18; B = 0.4; C = 0.4; Text = "al modificar conumo en litros o precio x litro, actualizar consumo `$ y diesel calculado (en el segundo fieldset)" },
    @{ Row = 19; B = 0.4; C = 0.4; Text = "al modificar kilometraje inicial o kilometraje final, actualizar kilometraje recorrido" },
    @{ Row = 20; B = 0.2; C = 0.2; Text = "implementar calculo de consumo en litros" },
    @{ Row = 21; B = 0.2; C = 0.2; Text = "implementar calculo de consumo en pesos" },
    @{ Row = 22; B = 0.2; C = 0.2; Text = "implementar calculo de kilometraje recorrido" },
    @{ Row = 23; B = 0.2; C = 0.2; Text = "calcular diesel lt" },
    @{ Row = 24; B = 0.2; C = 0.2; Text = "calcular diesel pesos" },
    @{ Row = 25; B = 0.2; C = 0.2; Text = "mostrar diferenca contra calculado" },
    @{ Row = 26; B = 0.2; C = 0.2; Text = "mostrar diferencia contra medido" }
)

# Numeric columns first, in natural row order.
foreach ($task in $tasks) {
    $r = $task.Row
    $ws2.Range("B$r").Value = $task.B
    $ws2.Range("C$r").Value = $task.C
}

# Text column second, but written in the exact order the author's shared
# string table shows them (row 19's text was registered before row 18's),
# so new <si> entries land at the same indices as the target workbook.
$textOrder = @(16, 17, 19, 18, 20, 21, 22, 23, 24, 25, 26)
foreach ($r in $textOrder) {
    $task = $tasks | Where-Object { $_.Row -eq $r }
    $ws2.Range("D$r").Value = $task.Text
}

# Stamp the D16:D26 cells with the same "task label" style (xfId 7: left
# aligned, indent 1) already used by D11/D13 etc., via copy/paste-special
# of formats only so no new style entries get minted.
$ws2.Range("D11").Copy()
$ws2.Range("D16:D26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- View state: restore selections on both sheets -------------------------
$ws1.Activate()
$ws1.Range("C10").Select()

$ws2.Activate()
$ws2.Range("D5").Select()
